$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17, shifting existing rows 17-128 down to 18-129.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with this week's record.
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "Terminal La Palmera de La Serena"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44635
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 100112040
$ws.Range("G17").Value = "Cilantro"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 1600
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 2750
$ws.Range("N17").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O17").Value = "Provincia del Elquí"
$ws.Range("P17").Value = 1833
$ws.Range("Q17").Value = 1.5
$ws.Range("R17").Value = "Hortaliza"

# Match the date-formatted style used by column D in the other rows.
$ws.Range("D17").NumberFormat = $ws.Range("D18").NumberFormat
